$d = $word.ActiveDocument

# --- Paragraph 1 formatting: add a (border-less style) paragraph border
#     with 5-twip spacing on all sides, and widen the left indent from
#     120 -> 225 twips (6pt -> 11.25pt). ---
$p1 = $d.Paragraphs(1)
$p1.Format.LeftIndent = 11.25
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5

# --- Paragraph 1 text: the first run becomes the new placeholder id and
#     the trailing " " run is removed entirely. ---
$oldId = "**ID__AFFARS_pgi_5307_topic_11__ID**"
$newId = "**ID__AFFARS_AFMC_PGI_5307_105__ID**"

# Delete the trailing space run first (from the tail end) so the first
# run's character offsets stay valid while we still need them.
$spaceRange = $d.Range($oldId.Length, $oldId.Length + 1)
$spaceRange.Text = ""

# Now replace the id text itself.
$idRange = $d.Range(0, $oldId.Length)
$idRange.Text = $newId
